$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.175.39"
$ws.Range("E2").Value = "  +1.37%  "

$ws.Range("D3").Value = "2.002.80"
$ws.Range("E3").Value = "  +2.09%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.94"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +2.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0806"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.32%  "

$ws.Range("D14").Value = "2.293.14"
$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.842"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.26%  "

$ws.Range("D17").Value = "2.011.28"
$ws.Range("E17").Value = "  +2.44%  "

$ws.Range("D18").Value = "37.102.26"
$ws.Range("E18").Value = "  +1.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.56%  "

$ws.Range("D20").Value = "0.0₃0863"
$ws.Range("E20").Value = "  +1.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.20"
$ws.Range("D22").Style = "Normal"

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.39%  "

$ws.Range("E27").Value = "  -4.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.63"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.50%  "

$ws.Range("E31").Value = "  +1.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0655"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.59%  "

$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("E37").Value = "  +2.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.99%  "

$ws.Range("E39").Value = "  -3.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0986"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.38%  "

$ws.Range("E41").Value = "  +0.88%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.59%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.91%  "

$ws.Range("D46").Value = "1.368.69"
$ws.Range("E46").Value = "  -0.21%  "

$ws.Range("E47").Value = "  +1.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +13.36%  "

$ws.Range("E50").Value = "  +5.34%  "

$ws.Range("E51").Value = "  -0.16%  "
